$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.629.99"
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = "'1.882.41"
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'249.39"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = "'0.4751"
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").Value = "'0.06534"
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'21.93"
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").Value = "'0.07750"
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = "'96.74"
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = "'0.7381"
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").Value = "'1.880.70"
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("D16").Value = "'274.46"
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").Value = "'30.602.69"
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("D19").Value = "'0.000007530"
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = "'2.127.24"
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = "'5.342"
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = "'6.239"
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("D26").Value = "'163.88"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = "'18.85"
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("D28").Value = "'1.907"
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("D30").Value = "'0.09706"
$ws.Range("E30").Value = '  -2.63%  '
$ws.Range("D31").Value = "'1.506"
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").Value = "'4.281"
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("D33").Value = "'4.151"
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("D34").Value = "'0.04866"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'0.6987"
$ws.Range("D38").Value = "'0.01905"
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").Value = "'2.794"
$ws.Range("E39").Value = '  +2.32%  '
$ws.Range("D40").Value = "'6.297"
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").Value = "'74.79"
$ws.Range("E41").Value = '  +6.17%  '
$ws.Range("E42").Value = '  +4.66%  '
$ws.Range("D43").Value = "'0.4248"
$ws.Range("E43").Value = '  +1.27%  '
$ws.Range("D44").Value = "'0.8405"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = "'102.58"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = "'9.359"
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").Value = "'7.046"
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").Value = "'35.60"
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").Value = "'915.94"
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("D51").Value = "'0.05756"
$ws.Range("E51").Value = '  +2.08%  '
